# edit.ps1 -- apply the changes described by the target diff:
#   1. Insert 4 new "Title and Content" slides at positions 14-17
#      (pushes the old Deep-Learning / Model-Evaluation slides down to 18-21).
#   2. Fill in title + body text for each of the 4 new slides.
#   3. Give the "CLUSTERING - RESULTS" slide's content placeholder an
#      explicit position/size (a:xfrm) instead of inheriting it from the
#      layout.

function EmuToPt($emu) {
    # PowerPoint COM exposes Left/Top/Width/Height in points while the
    # underlying XML stores EMU (1 pt = 12700 EMU). Nudge by half an EMU
    # before dividing so the round-trip lands back on the exact EMU value.
    return ($emu + 0.5) / 12700.0
}

$p = $ppt.ActivePresentation

# ------------------------------------------------------------------
# 1) Give "Content Placeholder 3" on the CLUSTERING - RESULTS slide
#    (slide 13) an explicit xfrm.
# ------------------------------------------------------------------
$clusterResults = $p.Slides.Item(13)
$contentPh = $clusterResults.Shapes.Item(2)
$contentPh.Left = EmuToPt(1141412)
$contentPh.Top = EmuToPt(2057399)
$contentPh.Width = EmuToPt(4876800)
$contentPh.Height = EmuToPt(3124201)

# ------------------------------------------------------------------
# 2) Insert four new "Title and Content" slides at index 14..17.
#    Layout 2 == the "Title and Content" custom layout (slideLayout2.xml).
# ------------------------------------------------------------------
$lsq = [char]0x2018
$rsq = [char]0x2019

$slideUserUser = $p.Slides.Add(14, 2)
$slideUserUser.Shapes.Item(1).TextFrame.TextRange.Text = "User-user based collaborative filtering"
$bodyUserUser = $slideUserUser.Shapes.Item(2).TextFrame.TextRange
$bodyUserUser.Text = "Created Utility matrix for users vs business filled with rating `r" + `
    "Centralized each rating vector of user to mean zero`r" + `
    "Replace the missing values with zero`r" + `
    "Averaged rating value would be 0, positive or negative where 0 represents mean rating or missing ratings and positive represents higher ratings than mean rating"

$slideRatingPred = $p.Slides.Add(15, 2)
$slideRatingPred.Shapes.Item(1).TextFrame.TextRange.Text = "Rating prediction"
$bodyRatingPred = $slideRatingPred.Shapes.Item(2).TextFrame.TextRange
$bodyRatingPred.Text = "Predicted the ratings that user will give to the new items`r" + `
    "Find top k similar users which have already rated the item " + $lsq + "I" + $rsq + " and returned the weighted averaged of ratings with weight as similarity"

$slideTop20 = $p.Slides.Add(16, 2)
$slideTop20.Shapes.Item(1).TextFrame.TextRange.Text = "Recommending top 20 business to user"
$bodyTop20 = $slideTop20.Shapes.Item(2).TextFrame.TextRange
$bodyTop20.Text = "Predict rating for business for the user which one have not used yet`r" + `
    "Sort the array on basis of rating`r" + `
    "Return the top 20 business"

$slideResult = $p.Slides.Add(17, 2)
$slideResult.Shapes.Item(1).TextFrame.TextRange.Text = "RESULT"

Write-Host "Slides after edit: $($p.Slides.Count)"
